# D4ASCopyright.docx update
# --------------------------
# Bumps the two version/year digits called out in the commit:
#   "Dynamo for Advance Steel 2025" -> "Dynamo for Advance Steel 2026"
#   "(c) 2024 Autodesk, Inc."       -> "(c) 2025 Autodesk, Inc."
#
# In the source document both years already live in their own tiny runs
# (just the trailing digit, e.g. <w:t>5</w:t> / <w:t>4</w:t>), so rather
# than replacing the whole "Advance Steel 2025" / "2024 Autodesk, Inc."
# phrase (which would reflow/merge a lot of unrelated runs), we locate
# each phrase with Find and then overwrite only the single trailing
# digit that actually needs to change.

$d = $word.ActiveDocument

# ---- 1) Title line: "...Advance Steel 2025" -> "...Advance Steel 2026"
$rng1 = $d.Content.Duplicate
$found1 = $rng1.Find.Execute("Steel 2025", $false, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'Steel 2025' in the document."
}
$digit1 = $d.Range($rng1.End - 1, $rng1.End)
$digit1.Text = "6"

# ---- 2) Copyright line: "(c) 2024 Autodesk, Inc." -> "(c) 2025 Autodesk, Inc."
$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute("2024 Autodesk, Inc.", $false, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find '2024 Autodesk, Inc.' in the document."
}
$yearEnd = $rng2.Start + 4   # length of "2024"
$digit2 = $d.Range($yearEnd - 1, $yearEnd)
$digit2.Text = "5"

$d.Save()
